$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '61.557.12'

# Row 3
$ws.Range("D3").Value = '2.876.05'
$ws.Range("E3").Value = '  -2.53%  '

# Row 4
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '564.73'
$ws.Range("E5").Value = '  -4.42%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.27'
$ws.Range("E6").Value = '  -3.28%  '

# Row 7
$ws.Range("E7").Value = '  +0.01%  '

# Row 8
$ws.Range("E8").Value = '  -1.31%  '

# Row 9
$ws.Range("D9").Value = '2.877.25'
$ws.Range("E9").Value = '  -2.44%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.87'
$ws.Range("E10").Value = '  -2.10%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.146'
$ws.Range("E11").Value = '  -1.62%  '

# Row 12
$ws.Range("E12").Value = '  -2.02%  '

# Row 13
$ws.Range("E13").Value = '  -1.27%  '

# Row 14
$ws.Range("E14").Value = '  -2.28%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.126'
$ws.Range("E15").Value = '  +0.18%  '

# Row 16
$ws.Range("D16").Value = '3.354.29'
$ws.Range("E16").Value = '  -2.46%  '

# Row 17
$ws.Range("D17").Value = '61.543.27'
$ws.Range("E17").Value = '  -2.06%  '

# Row 18
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.51'
$ws.Range("E18").Value = '  -2.54%  '

# Row 19
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '2.856.56'
$ws.Range("E19").Value = '  -3.09%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '429.38'
$ws.Range("E20").Value = '  -1.87%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.98'
$ws.Range("E21").Value = '  -3.38%  '

# Row 22
$ws.Range("E22").Value = '  -2.23%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.80'
$ws.Range("E23").Value = '  -2.41%  '

# Row 24
$ws.Range("E24").Value = '  -2.13%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.91'
$ws.Range("E25").Value = '  +0.98%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.11'
$ws.Range("E26").Value = '  -10.26%  '

# Row 27
$ws.Range("E27").Value = '  +0.10%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.02'
$ws.Range("E28").Value = '  -4.52%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000108'
$ws.Range("E29").Value = '  +7.45%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.01'
$ws.Range("E30").Value = '  -3.60%  '

# Row 31
$ws.Range("E31").Value = '  -4.37%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.01'
$ws.Range("E32").Value = '  -8.56%  '

# Row 33
$ws.Range("E33").Value = '  +0.06%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.106'
$ws.Range("E34").Value = '  -2.37%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '25.53'
$ws.Range("E35").Value = '  -2.88%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.949'
$ws.Range("E36").Value = '  -4.55%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.34'
$ws.Range("E37").Value = '  -4.44%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '48.82'
$ws.Range("E38").Value = '  -1.80%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.76'
$ws.Range("E39").Value = '  -8.03%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.89'
$ws.Range("E40").Value = '  -5.69%  '

# Row 41
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.114'
$ws.Range("E41").Value = '  -1.17%  '

# Row 42
$ws.Range("B42").Value = 'Cosmos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.12'
$ws.Range("E42").Value = '  -3.23%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.38'
$ws.Range("E43").Value = '  -1.90%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.268'
$ws.Range("E44").Value = '  -3.01%  '

# Row 45
$ws.Range("D45").Value = '2.677.76'
$ws.Range("E45").Value = '  -0.15%  '

# Row 46
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0335'
$ws.Range("E46").Value = '  -0.47%  '

# Row 47
$ws.Range("B47").Value = 'Monero'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '130.62'
$ws.Range("E47").Value = '  -3.39%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '344.58'
$ws.Range("E48").Value = '  -3.17%  '

# Row 49
$ws.Range("E49").Value = '  -0.04%  '

# Row 50
$ws.Range("E50").Value = '  -1.76%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.36'
$ws.Range("E51").Value = '  -5.48%  '
